$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value = 3544.6
$ws.Range("I69").Value = 3563
$ws.Range("J69").Value = 3540
$ws.Range("K69").Value = 10689
$ws.Range("L69").Value = 10620
$ws.Range("M69").Value = -9815
$ws.Range("N69").Value = -12368

# Row 72
$ws.Range("H72").Value = 3544.6
$ws.Range("I72").Value = 3563
$ws.Range("J72").Value = 3540
$ws.Range("K72").Value = 32067
$ws.Range("L72").Value = 31860
$ws.Range("M72").Value = -27699
$ws.Range("N72").Value = -40596

# Row 100
$ws.Range("H100").Value = 2866.842
$ws.Range("I100").Value = 2234
$ws.Range("K100").Value = 2234
$ws.Range("M100").Value = -1693

# Row 138
$ws.Range("H138").Value = 5248.1147
$ws.Range("I138").Value = 3979.4
$ws.Range("J138").Value = 5361.393
$ws.Range("K138").Value = 11938.2
$ws.Range("L138").Value = 16084.179
$ws.Range("M138").Value = -6798.200000000001
$ws.Range("N138").Value = -26364.179

# Row 141
$ws.Range("H141").Value = 254588.52
$ws.Range("I141").Value = 984.5897
$ws.Range("J141").Value = 1490907.6
$ws.Range("K141").Value = 2953.7691
$ws.Range("L141").Value = 4472722.800000001
$ws.Range("M141").Value = 2226.2309
$ws.Range("N141").Value = -4483082.800000001

$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

# Row 45
$ws.Range("H45").Value = 1357.0605
$ws.Range("I45").Value = 1019.7407
$ws.Range("J45").Value = 2875
$ws.Range("K45").Value = 1019.7407
$ws.Range("L45").Value = 2875
$ws.Range("M45").Value = -642.7406999999999
$ws.Range("N45").Value = -3629

# Row 61
$ws.Range("H61").Value = 3034.4243
$ws.Range("I61").Value = 1118
$ws.Range("J61").Value = 4446.5264
$ws.Range("K61").Value = 1118
$ws.Range("L61").Value = 4446.5264
$ws.Range("M61").Value = -906
$ws.Range("N61").Value = -4870.5264

# Row 74
$ws.Range("H74").Value = 651
$ws.Range("I74").Value = 590.7083
$ws.Range("J74").Value = 1133.3334
$ws.Range("K74").Value = 590.7083
$ws.Range("L74").Value = 1133.3334
$ws.Range("M74").Value = 283.2917
$ws.Range("N74").Value = -2881.3334

# Row 77
$ws.Range("H77").Value = 651
$ws.Range("I77").Value = 590.7083
$ws.Range("J77").Value = 1133.3334
$ws.Range("K77").Value = 2953.5415
$ws.Range("L77").Value = 5666.666999999999
$ws.Range("M77").Value = 1414.4585
$ws.Range("N77").Value = -14402.667

# Row 132
$ws.Range("H132").Value = 18520556
$ws.Range("I132").Value = 23810790
$ws.Range("J132").Value = 4735.5
$ws.Range("K132").Value = 71432370
$ws.Range("L132").Value = 14206.5
$ws.Range("M132").Value = -71429840
$ws.Range("N132").Value = -19266.5

# Row 136
$ws.Range("H136").Value = 3034.4243
$ws.Range("I136").Value = 1118
$ws.Range("J136").Value = 4446.5264
$ws.Range("K136").Value = 3354
$ws.Range("L136").Value = 13339.5792
$ws.Range("M136").Value = -804
$ws.Range("N136").Value = -18439.5792

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2184.7778
$ws.Range("I134").Value = 1743.68
$ws.Range("J134").Value = 3187.2727
$ws.Range("K134").Value = 5231.04
$ws.Range("L134").Value = 9561.8181
$ws.Range("M134").Value = -2696.04
$ws.Range("N134").Value = -14631.8181

# Row 141
$ws.Range("H141").Value = 37045.453
$ws.Range("J141").Value = 29642.857
$ws.Range("L141").Value = 29642.857
$ws.Range("N141").Value = -40002.857

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4238.4
$ws.Range("I31").Value = 2288
$ws.Range("J31").Value = 6351.3335
$ws.Range("K31").Value = 2288
$ws.Range("L31").Value = 6351.3335
$ws.Range("M31").Value = -1993
$ws.Range("N31").Value = -6941.3335

# Row 34
$ws.Range("H34").Value = 4238.4
$ws.Range("I34").Value = 2288
$ws.Range("J34").Value = 6351.3335
$ws.Range("K34").Value = 2288
$ws.Range("L34").Value = 6351.3335
$ws.Range("M34").Value = -2086
$ws.Range("N34").Value = -6755.3335

# Row 99
$ws.Range("H99").Value = 2785.5715
$ws.Range("I99").Value = 1249.75
$ws.Range("J99").Value = 4833.3335
$ws.Range("K99").Value = 1249.75
$ws.Range("L99").Value = 4833.3335
$ws.Range("M99").Value = 248.25
$ws.Range("N99").Value = -7829.3335

# Row 126
$ws.Range("H126").Value = 2785.5715
$ws.Range("I126").Value = 1249.75
$ws.Range("J126").Value = 4833.3335
$ws.Range("K126").Value = 3749.25
$ws.Range("L126").Value = 14500.0005
$ws.Range("M126").Value = -1279.25
$ws.Range("N126").Value = -19440.0005

# Row 134
$ws.Range("H134").Value = 1650.7142
$ws.Range("I134").Value = 748.1053000000001
$ws.Range("J134").Value = 3556.2222
$ws.Range("K134").Value = 2244.3159
$ws.Range("L134").Value = 10668.6666
$ws.Range("M134").Value = 290.6840999999999
$ws.Range("N134").Value = -15738.6666

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1371.9333
$ws.Range("I131").Value = 1871.25
$ws.Range("J131").Value = 1190.3636
$ws.Range("K131").Value = 5613.75
$ws.Range("L131").Value = 3571.0908
$ws.Range("M131").Value = -573.75
$ws.Range("N131").Value = -13651.0908

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2260
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 3310
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 3310
$ws.Range("M40").Value = -864
$ws.Range("N40").Value = -3582

# Row 132
$ws.Range("H132").Value = 2393.7437
$ws.Range("I132").Value = 1086.12
$ws.Range("J132").Value = 4728.7856
$ws.Range("K132").Value = 3258.36
$ws.Range("L132").Value = 14186.3568
$ws.Range("M132").Value = -728.3599999999997
$ws.Range("N132").Value = -19246.3568

# Row 136
$ws.Range("H136").Value = 1220.3684
$ws.Range("I136").Value = 841.4666999999999
$ws.Range("K136").Value = 2524.4001
$ws.Range("M136").Value = 25.59990000000016

# Row 137
$ws.Range("H137").Value = 29615
$ws.Range("J137").Value = 29615
$ws.Range("L137").Value = 29615
$ws.Range("N137").Value = -39815

# Row 138
$ws.Range("H138").Value = 30000
$ws.Range("J138").Value = 30000
$ws.Range("L138").Value = 30000
$ws.Range("N138").Value = -40280

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 825.55554
$ws.Range("I136").Value = 465.1111
$ws.Range("K136").Value = 1395.3333
$ws.Range("M136").Value = 1154.6667
